$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.901.22"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "1.878.25"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.74%  "
$ws.Range("D5").Value = "'325.21"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D7").Value = "'0.4596"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'0.3879"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("D9").Value = "'0.07868"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "'0.9852"
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("D11").Value = "'21.78"
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").Value = "1.906.62"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "'6.984"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").Value = "'5.647"
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("D15").Value = "'0.06962"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").Value = "'87.99"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").Value = "'0.000009977"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  -0.55%  "
$ws.Range("D21").Value = "28.921.42"
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").Value = "'2.100"
$ws.Range("E24").Value = "  +1.83%  "
$ws.Range("D25").Value = "'156.24"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").Value = "'19.26"
$ws.Range("D27").Value = "'6.038"
$ws.Range("E27").Value = "  +3.83%  "
$ws.Range("D28").Value = "'1.927"
$ws.Range("E28").Value = "  -1.47%  "
$ws.Range("D29").Value = "'117.20"
$ws.Range("E29").Value = "  -0.80%  "
$ws.Range("D30").Value = "'0.09337"
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").Value = "'0.9028"
$ws.Range("E31").Value = "  -2.11%  "
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("D33").Value = "'1.318"
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("D35").Value = "'1.180"
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("D36").Value = "'0.05756"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").Value = "'7.682"
$ws.Range("E39").Value = "  -1.43%  "
$ws.Range("D40").Value = "'0.5647"
$ws.Range("E40").Value = "  -0.27%  "
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("D42").Value = "'9.668"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("D43").Value = "'2.264"
$ws.Range("E43").Value = "  +3.83%  "
$ws.Range("D44").Value = "'11.84"
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("D45").Value = "'0.5344"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "'0.07043"
$ws.Range("E46").Value = "  -1.49%  "
$ws.Range("D47").Value = "'1.846"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D48").Value = "'112.90"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").Value = "'2.518"
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("E50").Value = "  -5.13%  "
$ws.Range("D51").Value = "'70.68"
$ws.Range("E51").Value = "  -0.11%  "
